$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 399
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

$ws.Range("B5").Value = 0.9975000000000001
$ws.Range("C5").Value = 0.0025
$ws.Range("D5").Value = 0
